$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the sheet to reflect the new "through" date
$ws.Name = "Through 2022-07-31"

# Update the header label in I1 (shared string "2022 (through 07-30)" -> "2022 (through 07-31)")
$ws.Range("I1").Value = "2022 (through 07-31)"

# Update the data values for 2022-08-08
$ws.Range("I8").Value = 169
$ws.Range("I14").Value = 975
